$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'58.307.42"
$ws.Range("E2").Value = "'  -1.48%  "
$ws.Range("D3").Value = "'2.298.09"
$ws.Range("E3").Value = "'  -0.94%  "
$ws.Range("E4").Value = "'  +0.10%  "
$ws.Range("D5").Value = "'534.10"
$ws.Range("E5").Value = "'  -3.54%  "
$ws.Range("D6").Value = "'131.60"
$ws.Range("E6").Value = "'  +0.11%  "
$ws.Range("E7").Value = "'  +0.05%  "
$ws.Range("D8").Value = "'0.586"
$ws.Range("E8").Value = "'  +2.71%  "
$ws.Range("D9").Value = "'2.295.07"
$ws.Range("E9").Value = "'  -1.02%  "
$ws.Range("D10").Value = "'0.0998"
$ws.Range("E10").Value = "'  -2.95%  "
$ws.Range("E11").Value = "'  -1.39%  "
$ws.Range("E12").Value = "'  -0.02%  "
$ws.Range("E13").Value = "'  -2.30%  "
$ws.Range("D14").Value = "'23.54"
$ws.Range("E14").Value = "'  -1.21%  "
$ws.Range("D15").Value = "'2.707.23"
$ws.Range("E15").Value = "'  -1.13%  "
$ws.Range("D16").Value = "'58.228.25"
$ws.Range("E16").Value = "'  -1.56%  "
$ws.Range("E17").Value = "'  -1.41%  "
$ws.Range("D18").Value = "'2.297.90"
$ws.Range("E18").Value = "'  -1.09%  "
$ws.Range("D19").Value = "'10.57"
$ws.Range("E19").Value = "'  -2.29%  "
$ws.Range("E20").Value = "'  -4.48%  "
$ws.Range("D21").Value = "'312.64"
$ws.Range("E21").Value = "'  -1.32%  "
$ws.Range("E22").Value = "'  -1.78%  "
$ws.Range("E23").Value = "'  -0.07%  "
$ws.Range("D24").Value = "'62.68"
$ws.Range("E24").Value = "'  -1.07%  "
$ws.Range("D25").Value = "'0.167"
$ws.Range("E25").Value = "'  -1.61%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "'  -0.13%  "
$ws.Range("D27").Value = "'8.03"
$ws.Range("E27").Value = "'  -3.58%  "
$ws.Range("D28").Value = "'1.26"
$ws.Range("E28").Value = "'  -5.69%  "
$ws.Range("D29").Value = "'170.55"
$ws.Range("E29").Value = "'  +0.21%  "
$ws.Range("E30").Value = "'  -3.94%  "
$ws.Range("E31").Value = "'  -2.11%  "
$ws.Range("D32").Value = "'5.75"
$ws.Range("E32").Value = "'  -2.15%  "
$ws.Range("E33").Value = "'  -3.01%  "
$ws.Range("E34").Value = "'  -3.55%  "
$ws.Range("E35").Value = "'  +0.01%  "
$ws.Range("E37").Value = "'  +0.10%  "
$ws.Range("E38").Value = "'  -4.01%  "
$ws.Range("D39").Value = "'3.88"
$ws.Range("E39").Value = "'  -3.37%  "
$ws.Range("D40").Value = "'38.53"
$ws.Range("E40").Value = "'  +0.05%  "
$ws.Range("E41").Value = "'  -4.08%  "
$ws.Range("D42").Value = "'140.74"
$ws.Range("E42").Value = "'  -1.94%  "
$ws.Range("D43").Value = "'289.15"
$ws.Range("E44").Value = "'  -0.92%  "
$ws.Range("E45").Value = "'  -0.10%  "
$ws.Range("E46").Value = "'  -1.24%  "
$ws.Range("D47").Value = "'0.556"
$ws.Range("E47").Value = "'  -0.43%  "
$ws.Range("D48").Value = "'18.09"
$ws.Range("E48").Value = "'  -3.15%  "
$ws.Range("E49").Value = "'  -2.25%  "
$ws.Range("D50").Value = "'10.94"
$ws.Range("E50").Value = "'  -0.94%  "
$ws.Range("E51").Value = "'  -0.64%  "
